$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename sheet
$ws.Name = "Sheet1"

# Insert 4 new rows at the top, pushing the existing table down
$ws.Rows("1:4").Insert()

# New header/title block content
$ws.Range("A1").Value = "MASTER PACKAGE"
$ws.Range("A2").Value = "WesternGlove Centric8 PROD"
$ws.Range("B2").Value = "M12225BVS563:KONRAD"
$ws.Range("C2").Value = "BOM"
$ws.Range("D2").Value = "MASTER"
$ws.Range("A3").Value = "Placements"

# Apply the same bordered/wrapped formatting used by the rest of the table
# (copy format from the row that used to be row 1, now row 5)
$ws.Range("A5:D5").Copy()
$ws.Range("A1:D3").PasteSpecial(-4122)

# Widen column B to fit the new, longer placements text
$ws.Columns("B").ColumnWidth = 21.6
